# Commited changes for Fluent wait and RowNumber in data drivern
#
# 1. Insert a new first column "RowNumber" (pushes all existing columns
#    one place to the right) and seed its single data value with 0.
# 2. Split the old "Miles" pair of 50/50 values across four columns
#    (DryVan / Refrigerated / Flatbed / Intermodal) as 25/25/25/25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column A ("RowNumber") - shifts every other column
# (and their data) one position to the right automatically.
$ws.Range("A1").EntireColumn.Insert()

$ws.Range("A1").Value = "RowNumber"
$ws.Range("A2").Value = 0

# The "Miles" values that used to live in two columns (50/50) are now
# spread across four columns (DryVan/Refrigerated/Flatbed/Intermodal).
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 25
$ws.Range("O2").Value = 25
$ws.Range("P2").Value = 25
